$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally contained 18 data rows, where the 5-row block of
# guide-RNA results (rows 1-5) repeated three times, interleaved with an
# extra "orphan" guide (chr11:59318773-59318796, no upstream match) that
# appeared twice in a row (rows 6-7) and once more at the end (row 18).
# The fix removes that orphan guide entirely, leaving just the clean
# 3x repeat of the 5-row block (18 - 3 = 15 rows).
#
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
